# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
# Adds a new "Number of employees / Assets / Turnover" breakdown table
# (rows 23-27) to the Malta MSME summary sheet, and relocates the
# "SME Performance Review EU" source citation from rows 26-27 down to
# rows 32-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, grab the text that currently lives in A26 / A27 so we can
# re-home it lower on the sheet (rows 32/33) without disturbing the
# existing shared-string entries.
$sourceLabel = $ws.Range("A26").Value2
$sourceText  = $ws.Range("A27").Value2

# New header row (row 23): Number of employees / Assets / Turnover
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B23:D23").Style = "title"

# Row 24: Micro / <10 / (blank) / (blank)
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").Value = "'"
$ws.Range("D24").Value = "'"

# Row 25: Small / <50 / (blank) / (blank)
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").Value = "'"
$ws.Range("D25").Value = "'"

# Row 26: Medium / <250 / (blank) / (blank)
$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "<250"
$ws.Range("C26").Value = "'"
$ws.Range("D26").Value = "'"

# Row 27: Large / >249 / (blank) / (blank)
$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">249"
$ws.Range("C27").Value = "'"
$ws.Range("D27").Value = "'"

$ws.Range("A24:D27").Style = "Normal"

# Re-home the source citation down to rows 32/33.
$ws.Range("A32").Value = $sourceLabel
$ws.Range("A33").Value = $sourceText
$ws.Range("A32").Style = "title"
$ws.Range("A33").Style = "source"
